# Scheduled-runner refresh of market-price-derived columns (H:N) on the
# Behemoth_Profits "Leve" sheets (ALC, ARM, BSM, CRP, GSM, LTW, WVR).
# Columns H..N are plain cached numbers (no formulas in this workbook),
# so the refreshed Universalis price snapshot is written back cell-by-cell.
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
# Row 5: Met a Sticky End / Animal Glue
$ws.Range("H5").Value = 461.4
$ws.Range("I5").Value = 76.75
$ws.Range("J5").Value = 2000
$ws.Range("K5").Value = 76.75
$ws.Range("L5").Value = 2000
$ws.Range("M5").Value = 38.25
$ws.Range("N5").Value = -2230

# Row 80: Cleansing the Wicked Humours / Hallowed Water
$ws.Range("H80").Value = 1589
$ws.Range("I80").Value = 749
$ws.Range("K80").Value = 2247
$ws.Range("M80").Value = -1249

# Row 83: Washing Away the Sins (L) / Hallowed Water
$ws.Range("H83").Value = 1589
$ws.Range("I83").Value = 749
$ws.Range("K83").Value = 6741
$ws.Range("M83").Value = -1749

# Row 100: Asking for a Friend / Beetle Glue
$ws.Range("H100").Value = 3117.8235
$ws.Range("I100").Value = 1335
$ws.Range("J100").Value = 3499.8572
$ws.Range("K100").Value = 1335
$ws.Range("L100").Value = 3499.8572
$ws.Range("M100").Value = -794
$ws.Range("N100").Value = -4581.8572

# Row 132: Fast-forwarding Flora / Growth Formula Lambda
$ws.Range("H132").Value = 9999.5
$ws.Range("I132").Value = 9999
$ws.Range("J132").Value = 10000
$ws.Range("K132").Value = 29997
$ws.Range("L132").Value = 30000
$ws.Range("M132").Value = -27467
$ws.Range("N132").Value = -35060

# Row 138: All-night Crafting / Cunning Craftsman's Tisane
$ws.Range("H138").Value = 2617.169
$ws.Range("I138").Value = 1795.5333
$ws.Range("J138").Value = 2815.9517
$ws.Range("K138").Value = 5386.5999
$ws.Range("L138").Value = 8447.8551
$ws.Range("M138").Value = -246.5999000000002
$ws.Range("N138").Value = -18727.8551

# Row 141: Remedy for Reason / Grade 1 Gemdraught of Mind
$ws.Range("H141").Value = 4437.88
$ws.Range("I141").Value = 4974.409
$ws.Range("K141").Value = 14923.227
$ws.Range("M141").Value = -9743.226999999999

$ws = $wb.Worksheets.Item("ARM")
# Row 61: Dealing with the Tough Stuff / Cobalt Ingot
$ws.Range("H61").Value = 19743966
$ws.Range("I61").Value = 15157584
$ws.Range("J61").Value = 50014090
$ws.Range("K61").Value = 15157584
$ws.Range("L61").Value = 50014090
$ws.Range("M61").Value = -15157372
$ws.Range("N61").Value = -50014514

# Row 74: As the Bolt Flies / Titanium Nugget
$ws.Range("H74").Value = 15301735
$ws.Range("I74").Value = 25005862
$ws.Range("J74").Value = 1438696.6
$ws.Range("K74").Value = 25005862
$ws.Range("L74").Value = 1438696.6
$ws.Range("M74").Value = -25004988
$ws.Range("N74").Value = -1440444.6

# Row 77: Heavy Metal Banned (L) / Titanium Nugget
$ws.Range("H77").Value = 15301735
$ws.Range("I77").Value = 25005862
$ws.Range("J77").Value = 1438696.6
$ws.Range("K77").Value = 125029310
$ws.Range("L77").Value = 7193483
$ws.Range("M77").Value = -125024942
$ws.Range("N77").Value = -7202219

# Row 110: Scheduled Maintenance / Deepgold Ingot
$ws.Range("H110").Value = 1866.5625
$ws.Range("I110").Value = 1824.3334
$ws.Range("K110").Value = 1824.3334
$ws.Range("M110").Value = 220.6666

# Row 132: Don't Bore Me, Ore Me / Mountain Chromite Ingot
$ws.Range("H132").Value = 418.8
$ws.Range("I132").Value = 418.8
$ws.Range("J132").Value = 0
$ws.Range("K132").Value = 1256.4
$ws.Range("L132").Value = 0
$ws.Range("M132").Value = 1273.6
$ws.Range("N132").ClearContents()

# Row 136: Metal with Mettle / Cobalt Tungsten Ingot
$ws.Range("H136").Value = 19743966
$ws.Range("I136").Value = 15157584
$ws.Range("J136").Value = 50014090
$ws.Range("K136").Value = 45472752
$ws.Range("L136").Value = 150042270
$ws.Range("M136").Value = -45470202
$ws.Range("N136").Value = -150047370

$ws = $wb.Worksheets.Item("BSM")
# Row 105: Ingot to Wing It / Molybdenum Ingot
$ws.Range("H105").Value = 2390.8262
$ws.Range("I105").Value = 2162
$ws.Range("J105").Value = 2490.9375
$ws.Range("K105").Value = 2162
$ws.Range("L105").Value = 2490.9375
$ws.Range("M105").Value = -415
$ws.Range("N105").Value = -5984.9375

# Row 134: Ruthenium Supremium / Ruthenium Ingot
$ws.Range("H134").Value = 527650.44
$ws.Range("I134").Value = 1343.7646
$ws.Range("K134").Value = 4031.2938
$ws.Range("M134").Value = -1496.2938

$ws = $wb.Worksheets.Item("CRP")
# Row 58: You Do the Heavy Lifting / Mahogany Lumber
$ws.Range("H58").Value = 2389.4546
$ws.Range("I58").Value = 1297.7407
$ws.Range("J58").Value = 7302.1665
$ws.Range("K58").Value = 1297.7407
$ws.Range("L58").Value = 7302.1665
$ws.Range("M58").Value = -1094.7407
$ws.Range("N58").Value = -7708.1665

# Row 119: Off to a Good Staff / Dwarven Lignum Cane
$ws.Range("H119").Value = 0
$ws.Range("J119").Value = 0
$ws.Range("L119").Value = 0
$ws.Range("N119").ClearContents()

# Row 134: Wood You Be Quiet / Ceiba Lumber
$ws.Range("H134").Value = 2646.8276
$ws.Range("I134").Value = 1288.3
$ws.Range("K134").Value = 3864.9
$ws.Range("M134").Value = -1329.9

# Row 136: Turali Quality / Dark Mahogany Lumber
$ws.Range("H136").Value = 2389.4546
$ws.Range("I136").Value = 1297.7407
$ws.Range("J136").Value = 7302.1665
$ws.Range("K136").Value = 3893.2221
$ws.Range("L136").Value = 21906.4995
$ws.Range("M136").Value = -1343.2221
$ws.Range("N136").Value = -27006.4995

$ws = $wb.Worksheets.Item("GSM")
# Row 126: Gold Rush Order / Phrygian Gold Ingot
$ws.Range("H126").Value = 5111.4287
$ws.Range("I126").Value = 780
$ws.Range("J126").Value = 5833.3335
$ws.Range("K126").Value = 2340
$ws.Range("L126").Value = 17500.0005
$ws.Range("M126").Value = 130
$ws.Range("N126").Value = -22440.0005

$ws = $wb.Worksheets.Item("LTW")
# Row 40: Best Served Toad / Toad Leather
$ws.Range("H40").Value = 2636.1
$ws.Range("I40").Value = 1849.3043
$ws.Range("J40").Value = 5221.2856
$ws.Range("K40").Value = 1849.3043
$ws.Range("L40").Value = 5221.2856
$ws.Range("M40").Value = -1713.3043
$ws.Range("N40").Value = -5493.2856

# Row 68: You Could Say It's a Moving Target / Wyvern Leather
$ws.Range("H68").Value = 4981.125
$ws.Range("I68").Value = 5179.6
$ws.Range("J68").Value = 4890.909
$ws.Range("K68").Value = 5179.6
$ws.Range("L68").Value = 4890.909
$ws.Range("M68").Value = -4430.6
$ws.Range("N68").Value = -6388.909

# Row 71: They Call It Bloody Mary (L) / Wyvern Leather
$ws.Range("H71").Value = 4981.125
$ws.Range("I71").Value = 5179.6
$ws.Range("J71").Value = 4890.909
$ws.Range("K71").Value = 25898
$ws.Range("L71").Value = 24454.545
$ws.Range("M71").Value = -22154
$ws.Range("N71").Value = -31942.545

# Row 121: A Shoe In / Swallowskin Shoes of Healing
$ws.Range("H121").Value = 0
$ws.Range("J121").Value = 0
$ws.Range("L121").Value = 0
$ws.Range("N121").ClearContents()

# Row 132: Tenets of Tanning / Silver Lobo Leather
$ws.Range("H132").Value = 638479
$ws.Range("I132").Value = 19876
$ws.Range("J132").Value = 1669484
$ws.Range("K132").Value = 59628
$ws.Range("L132").Value = 5008452
$ws.Range("M132").Value = -57098
$ws.Range("N132").Value = -5013512

# Row 136: Respect for Br'aax / Br'aax Leather
$ws.Range("H136").Value = 146464.67
$ws.Range("I136").Value = 37332.332
$ws.Range("K136").Value = 111996.996
$ws.Range("M136").Value = -109446.996

$ws = $wb.Worksheets.Item("WVR")
# Row 121: Healing Headwear / Dwarven Cotton Petasos of Healing
$ws.Range("H121").Value = 50000
$ws.Range("J121").Value = 50000
$ws.Range("L121").Value = 50000
$ws.Range("N121").Value = -53494

# Row 132: Comfy Cabins / Snow Cotton Cloth
$ws.Range("H132").Value = 5588.75
$ws.Range("I132").Value = 4455.4287
$ws.Range("K132").Value = 13366.2861
$ws.Range("M132").Value = -10836.2861
